$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 0.85
$ws.Range("F2").Value = 85.25

# Row 4
$ws.Range("E4").Value = 0.8
$ws.Range("F4").Value = 55.94

# Row 5
$ws.Range("E5").Value = 0.85
$ws.Range("F5").Value = 85.25

# Row 10
$ws.Range("E10").Value = 0.85
$ws.Range("F10").Value = 71.40000000000001
